$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the question text cells with HTML-formatted versions ---
# Order matters: the underlying shared-string table appends newly
# edited strings in the order the edits are made, and the target
# workbook expects them appended in this exact sequence (B4, B12,
# B19, B18, B20) so that the new entries land at shared-string
# indices 120-124 respectively.

$ws.Range("B4").Value = "Which of the following statements are TRUE about the .NET CLR?<br/>`n1.It provides a language-neutral development & execution environment.<br/>`n2.It ensures that an application would not be able to access memory that it is not authorized to access.<br/>`n3.It provides services to run ""managed"" applications.<br/>`n4.The resources are garbage collected.<br/>`n5.It provides services to run ""unmanaged"" applications."

$ws.Range("B12").Value = "Which of the statements are true ?<br/>`nI. Function overloading is done at compile time.<br/>`nII. Protected members are accessible to the member of derived class.<br/>`nIII. A derived class inherits constructors and destructors.<br/>`nIV. A friend function can be called like a normal function.<br/>`n<b>V. Nested class is a derived class.</b> <br/>"

$ws.Range("B19").Value = "<b>Which of the following statements are correct about constructors in C#.NET?</b><br/>`n1)Constructors cannot be overloaded.<br/>`n2)Constructors always have the name same as the name of the class.<br/>`n3)Constructors are never called explicitly.<br/>`n4)Constructors never return any value.<br/>`n5)Constructors allocate space for the object in memory.<br/>"

$ws.Range("B18").Value = "Which of the following is NOT an Arithmetic operator in C#.NET?<br/>`nA) **<br/>`nB) /<br/>`nC) +<br/>`nD) %<br/>`nE)~`n"

$ws.Range("B20").Value = "<span style=""color: #ff0000"">Which of the following statements is correct about constructors in C#.NET?</span>"

# --- Row heights recomputed by Excel's autofit for the wrapped
# question-text cells after the above edits ---
$ws.Rows(4).RowHeight = 147
$ws.Rows(12).RowHeight = 129.6
$ws.Rows(18).RowHeight = 115.2
$ws.Rows(19).RowHeight = 129.6

# --- Selection / active cell moved in the saved view state ---
$ws.Range("B20").Select()
